# Intro to ML completed
# Adds a computed "Average Sales House Prices" column (R) to Sheet1:
# header in R1, and for every data row 2-85, the arithmetic mean of the
# fifteen price columns C:Q, formatted with a plain integer (0) number
# format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column R.
$ws.Range("R1").Value = "Average Sales House Prices"

# Fill R2:R85 with the row average of C:Q, baked down to a static value
# (matching the workbook's existing style of literal values rather than
# live formulas).
for ($r = 2; $r -le 85; $r++) {
    $cell = $ws.Range("R$r")
    $cell.Formula = "=AVERAGE(C$r`:Q$r)"
    $cell.Value = $cell.Value2
    $cell.NumberFormat = "0"
}

# Column R width, sized to fit its contents.
$ws.Columns.Item(18).ColumnWidth = 11.83

# Match the saved selection/active cell.
$ws.Range("R2").Select()
